$d = $word.ActiveDocument
$d.Content.Find.Execute("https://wolof", $false, $false, $false, $false, $false, $true, 1, $false, "https://wolof", 2) | Out-Null
Write-Output "done"
